# Update base "Azerbaijan Premier League" - atualização de bases das ligas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix/replace existing row 147 (id = 145) with corrected match data
# ---------------------------------------------------------------------------
$ws.Range("B147").Value = 7011639
$ws.Range("E147").Value = 45395.39583333334
$ws.Range("F147").Value = "FK Gabala"
$ws.Range("G147").Value = "PFK Turan Tovuz"
$ws.Range("H147").Value = 1
$ws.Range("I147").Value = 2
$ws.Range("J147").Value = "A"
$ws.Range("K147").Value = 3
$ws.Range("L147").Value = 3.5
$ws.Range("M147").Value = 2
$ws.Range("N147").Value = 2.8
$ws.Range("O147").Value = 3.4
$ws.Range("P147").Value = 2.15
$ws.Range("Q147").Value = 0.25
$ws.Range("R147").Value = 1.85
$ws.Range("S147").Value = 1.95
$ws.Range("T147").Value = 2.25
$ws.Range("U147").Value = 1.75
$ws.Range("V147").Value = 1.95
$ws.Range("W147").Value = -1
$ws.Range("X147").Value = -1
$ws.Range("Y147").Value = 1.15
$ws.Range("Z147").Value = -1
$ws.Range("AA147").Value = 0.95
$ws.Range("AB147").Value = 0.75
$ws.Range("AC147").Value = -1

# ---------------------------------------------------------------------------
# 2) Append a brand-new row 148 (id = 146) for another match
# ---------------------------------------------------------------------------
$ws.Range("A148").Value = 146
$ws.Range("B148").Value = 7011637
$ws.Range("C148").Value = "Azerbaijan Premier League"
$ws.Range("D148").Value = "Azerbaijan Premier League"
$ws.Range("E148").Value = 45395.5
$ws.Range("F148").Value = "Neftchi Baku"
$ws.Range("G148").Value = "Araz FK"
$ws.Range("H148").Value = 3
$ws.Range("I148").Value = 0
$ws.Range("J148").Value = "H"
$ws.Range("K148").Value = 1.833
$ws.Range("L148").Value = 3.4
$ws.Range("M148").Value = 3.6
$ws.Range("N148").Value = 1.8
$ws.Range("O148").Value = 3.4
$ws.Range("P148").Value = 3.8
$ws.Range("Q148").Value = -0.5
$ws.Range("R148").Value = 1.8
$ws.Range("S148").Value = 2
$ws.Range("T148").Value = 2.25
$ws.Range("U148").Value = 1.8
$ws.Range("V148").Value = 2
$ws.Range("W148").Value = 0.8
$ws.Range("X148").Value = -1
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = 0.8
$ws.Range("AA148").Value = -1
$ws.Range("AB148").Value = 0.8
$ws.Range("AC148").Value = -1

# Replicate the cell formatting used by the rest of the table:
#  - column A uses the bold/centered/bordered "id" style
#  - column E uses the custom date-time number format
$ws.Range("A147").Copy()
$ws.Range("A148").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E147").Copy()
$ws.Range("E148").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "Update complete"
